# Updates generated data in "北京-漫展信息.xlsx":
# - "想去人数" (interest count) column F bumped for a number of rows
# - Cover image URL (column I) refreshed for one event row
# across sheets 展览 (1), 演出 (2), 本地生活 (3), 全部类型 (4).

$wb = $excel.ActiveWorkbook

$newCover = "//i0.hdslb.com/bfs/openplatform/202405/iR6rV5311717039317028.jpeg"

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 7757
$ws1.Range("F3").Value = 7757
$ws1.Range("F5").Value = 7916
$ws1.Range("F9").Value = 6733
$ws1.Range("F10").Value = 3400
$ws1.Range("F12").Value = 3738
$ws1.Range("F15").Value = 49
$ws1.Range("F16").Value = 72
$ws1.Range("F20").Value = 53
$ws1.Range("F21").Value = 328
$ws1.Range("F24").Value = 3885
$ws1.Range("F26").Value = 376
$ws1.Range("F27").Value = 958
$ws1.Range("F28").Value = 291
$ws1.Range("F29").Value = 1506
$ws1.Range("F33").Value = 1911
$ws1.Range("F35").Value = 53
$ws1.Range("F37").Value = 63
$ws1.Range("F38").Value = 3751
$ws1.Range("F39").Value = 337
$ws1.Range("F40").Value = 285
$ws1.Range("F43").Value = 555
$ws1.Range("F45").Value = 1447
$ws1.Range("F48").Value = 564
$ws1.Range("I49").Value = $newCover

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 418
$ws2.Range("F13").Value = 92
$ws2.Range("F17").Value = 152

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 141

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 141
$ws4.Range("F5").Value = 7757
$ws4.Range("F6").Value = 7757
$ws4.Range("F7").Value = 7916
$ws4.Range("F10").Value = 6733
$ws4.Range("F11").Value = 3400
$ws4.Range("F12").Value = 3738
$ws4.Range("F14").Value = 49
$ws4.Range("F15").Value = 72
$ws4.Range("F19").Value = 53
$ws4.Range("F20").Value = 328
$ws4.Range("F23").Value = 3885
$ws4.Range("F27").Value = 376
$ws4.Range("F28").Value = 958
$ws4.Range("F29").Value = 291
$ws4.Range("F30").Value = 1506
$ws4.Range("F34").Value = 1911
$ws4.Range("F36").Value = 53
$ws4.Range("F38").Value = 92
$ws4.Range("F39").Value = 3751
$ws4.Range("F40").Value = 337
$ws4.Range("F41").Value = 285
$ws4.Range("F44").Value = 555
$ws4.Range("F45").Value = 152
$ws4.Range("F46").Value = 1447
$ws4.Range("F49").Value = 564
$ws4.Range("I50").Value = $newCover
